# Update "想去人数" (F column) and a couple of "最低票价" (G column) values
# across the four sheets, per the source diff.

$wb = $excel.ActiveWorkbook

function Set-FValue {
    param($ws, [int]$row, [double]$value)
    $ws.Cells.Item($row, 6).Value = $value
}

function Set-GValue {
    param($ws, [int]$row, [double]$value)
    $ws.Cells.Item($row, 7).Value = $value
}

# Sheet: 展览
$wsExhibit = $wb.Worksheets.Item("展览")
Set-FValue $wsExhibit 2 590
Set-FValue $wsExhibit 3 267
Set-FValue $wsExhibit 5 745
Set-FValue $wsExhibit 6 391
Set-GValue $wsExhibit 6 36
Set-FValue $wsExhibit 8 165
Set-FValue $wsExhibit 10 234
Set-FValue $wsExhibit 11 6124
Set-FValue $wsExhibit 12 63
Set-FValue $wsExhibit 16 550
Set-FValue $wsExhibit 17 365
Set-FValue $wsExhibit 21 715
Set-FValue $wsExhibit 22 167
Set-FValue $wsExhibit 25 1028
Set-FValue $wsExhibit 26 66
Set-FValue $wsExhibit 27 1853
Set-FValue $wsExhibit 28 508

# Sheet: 演出
$wsShow = $wb.Worksheets.Item("演出")
Set-FValue $wsShow 3 273
Set-FValue $wsShow 4 53
Set-FValue $wsShow 5 273
Set-FValue $wsShow 6 302

# Sheet: 本地生活
$wsLocal = $wb.Worksheets.Item("本地生活")
Set-FValue $wsLocal 2 267

# Sheet: 全部类型
$wsAll = $wb.Worksheets.Item("全部类型")
Set-FValue $wsAll 2 267
Set-FValue $wsAll 3 590
Set-FValue $wsAll 4 267
Set-FValue $wsAll 6 745
Set-FValue $wsAll 8 391
Set-GValue $wsAll 8 36
Set-FValue $wsAll 10 165
Set-FValue $wsAll 12 234
Set-FValue $wsAll 13 6124
Set-FValue $wsAll 14 63
Set-FValue $wsAll 16 273
Set-FValue $wsAll 19 550
Set-FValue $wsAll 20 365
Set-FValue $wsAll 22 53
Set-FValue $wsAll 25 273
Set-FValue $wsAll 26 302
Set-FValue $wsAll 28 715
Set-FValue $wsAll 32 167
Set-FValue $wsAll 35 1028
Set-FValue $wsAll 36 66
Set-FValue $wsAll 37 1853
Set-FValue $wsAll 38 508
